# [LocTT] update 3 field export Main2
# Add 2 new template fields to the tab1 header row (row 18):
#   C18 -> {tab1.tfluid}
#   E18 -> {tab1.wc}
# (D18 already holds {tab1.pressure})

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E18 was previously blank with the default style; give it the same
# formatting as the rest of the header row (copy from D18) before writing
# the new placeholder text into it.
$ws.Range("D18").Copy()
$ws.Range("E18").PasteSpecial(-4122)

$ws.Range("C18").Value = "{tab1.tfluid}"
$ws.Range("E18").Value = "{tab1.wc}"

$ws.Range("E18").Select()
